# Update "想去人数" (want-to-go count) figures in column F for the
# "展览" and "全部类型" worksheets, reflecting a refreshed data pull.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1574
$ws1.Range("F5").Value = 270
$ws1.Range("F7").Value = 1493
$ws1.Range("F8").Value = 10198
$ws1.Range("F14").Value = 7070
$ws1.Range("F17").Value = 34

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1574
$ws4.Range("F5").Value = 270
$ws4.Range("F8").Value = 1493
$ws4.Range("F11").Value = 10198
$ws4.Range("F17").Value = 7070
$ws4.Range("F20").Value = 34
